$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.527469038963318
$ws.Range("B1").Value = 1.765408515930176
$ws.Range("C1").Value = 2.234684944152832
$ws.Range("D1").Value = 3.565440654754639
$ws.Range("E1").Value = 3.410730361938477
